# This script applies the following change to the "lojadosomeacessorios" sheet:
#  - Insert two new leading columns: "data" and "loja"
#  - Shift the former columns (nome, modelo, preco, politica, full, tipo, link)
#    two columns to the right (C..I)
#  - Refresh the product rows with updated data (5 rows instead of 6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new columns at the very left (A and B) -------------------
# This automatically shifts the existing A:G content (and its formatting)
# to C:I, carrying over header styles/borders etc.
$ws.Columns.Item(1).Insert()
$ws.Columns.Item(1).Insert()

# --- 2. Build the two new header cells, copying the header style from the
#        (now shifted) "nome" header cell C1 so the bold/border formatting
#        matches the rest of the header row ------------------------------
$ws.Range("C1").Copy($ws.Range("A1"))
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"

# --- 3. Write out the refreshed table data (header row stays the same for
#        the shifted columns; only the 9 columns below are rewritten) -----

# Row 2
$ws.Range("A2").Value = "30/07/2024"
$ws.Range("B2").Value = "lojadosomeacessorios"
$ws.Range("C2").Value = "Controle Longa Distância Jfa Acqua 1200 Resistente A Água"
$ws.Range("D2").Value = "ACQUA"
$ws.Range("E2").Value = 78.90000000000001
$ws.Range("F2").Value = "Baixo"
$ws.Range("G2").Value = "NA"
$ws.Range("H2").Value = "classico"
$ws.Range("I2").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27685629?pdp_filters=seller_id:296847653#searchVariation=MLB27685629&position=20&search_layout=grid&type=product&tracking_id=e9351b0a-97e9-41fd-9f30-557e972e8462"

# Row 3
$ws.Range("A3").Value = "30/07/2024"
$ws.Range("B3").Value = "lojadosomeacessorios"
$ws.Range("C3").Value = "Controle Longa Distância Jfa Acqua 1200 Resistente A Água"
$ws.Range("D3").Value = "ACQUA"
$ws.Range("E3").Value = 78.90000000000001
$ws.Range("F3").Value = "Baixo"
$ws.Range("G3").Value = "NA"
$ws.Range("H3").Value = "classico"
$ws.Range("I3").Value = "https://www.mercadolivre.com.br/controle-longa-distncia-jfa-acqua-1200-resistente-a-agua/p/MLB27687422?pdp_filters=seller_id:296847653#searchVariation=MLB27687422&position=23&search_layout=grid&type=product&tracking_id=e9351b0a-97e9-41fd-9f30-557e972e8462"

# Row 4
$ws.Range("A4").Value = "30/07/2024"
$ws.Range("B4").Value = "lojadosomeacessorios"
$ws.Range("C4").Value = "Controle Remoto Jfa Redline Wr Longo Alcance De Longa Distan"
$ws.Range("D4").Value = "CONTROLE WR"
$ws.Range("E4").Value = 121.47
$ws.Range("F4").Value = "Acima"
$ws.Range("G4").Value = "NA"
$ws.Range("H4").Value = "premium"
$ws.Range("I4").Value = "https://www.mercadolivre.com.br/controle-remoto-jfa-redline-wr-longo-alcance-de-longa-distan/p/MLB28557249?pdp_filters=seller_id:296847653#searchVariation=MLB28557249&position=29&search_layout=grid&type=product&tracking_id=e9351b0a-97e9-41fd-9f30-557e972e8462"

# Row 5
$ws.Range("A5").Value = "30/07/2024"
$ws.Range("B5").Value = "lojadosomeacessorios"
$ws.Range("C5").Value = "Amplificador 380w Jfa Ap380 Rms 4 Canais Crossover Lançament"
$ws.Range("D5").Value = "Sem Modelo"
$ws.Range("E5").Value = 351.4
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "NA"
$ws.Range("H5").Value = "classico"
$ws.Range("I5").Value = "https://produto.mercadolivre.com.br/MLB-3711906749-amplificador-380w-jfa-ap380-rms-4-canais-crossover-lancament-_JM#position%3D9%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D9f7769ff-b3be-47c3-9256-b9ecced16e39"

# Row 6 (replaces old rows 6 & 7 - "K1200" product is dropped, "Redline WR
# Longa Distância" product takes row 6 with refreshed price/link)
$ws.Range("A6").Value = "30/07/2024"
$ws.Range("B6").Value = "lojadosomeacessorios"
$ws.Range("C6").Value = "Controle Remoto Jfa Redline Wr Longa Distância 1200 Metros"
$ws.Range("D6").Value = "CONTROLE WR"
$ws.Range("E6").Value = 111.9
$ws.Range("F6").Value = "Acima"
$ws.Range("G6").Value = "NA"
$ws.Range("H6").Value = "classico"
$ws.Range("I6").Value = "https://produto.mercadolivre.com.br/MLB-2640116873-controle-remoto-jfa-redline-wr-longa-distncia-1200-metros-_JM#position%3D38%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D67704492-d5e4-4384-b133-1e38fcf0a70e"

# --- 4. Remove the now-obsolete old row 7 (original data had 6 rows, the
#        refreshed table only has 5) ---------------------------------------
$ws.Rows.Item(7).Delete()

# --- 5. Leave the selection on A1, matching the original sheet's selection -
[void]$ws.Range("A1").Select()
